# B6-PowerPoint.pptx edit
#
# 1) Three tables (on slides 14, 15 and 16) get their table style switched
#    from the custom "Table_0" style ({32EF5BB1-50FF-4347-BA12-31043699B090})
#    to the built-in style {B8E04126-19BD-49C2-9E0E-A682514183F2}.
#
# 2) The deck's theme colour scheme is swapped from the "Integral" / "Red
#    Violet" palette back to the standard Office palette (this is what the
#    author's diff shows as the two theme parts trading places - the theme
#    that is actually wired to the slide master/presentation ends up
#    carrying the plain "Office" colours).

$p = $ppt.ActivePresentation

# --- 1. Table styles -------------------------------------------------
$newStyleId = "{B8E04126-19BD-49C2-9E0E-A682514183F2}"

foreach ($slideIdx in 14, 15, 16) {
    $slide = $p.Slides.Item($slideIdx)
    for ($i = 1; $i -le $slide.Shapes.Count; $i++) {
        $shp = $slide.Shapes.Item($i)
        if ($shp.HasTable) {
            $shp.Table.ApplyStyle($newStyleId)
        }
    }
}

# --- 2. Theme colours --------------------------------------------------
$cs = $p.SlideMaster.ColorScheme

$cs.Colors(1).RGB  = 0          # dk1      -> 000000
$cs.Colors(2).RGB  = 16777215   # lt1      -> FFFFFF
$cs.Colors(3).RGB  = 6968388    # dk2      -> 44546A
$cs.Colors(4).RGB  = 15132391   # lt2      -> E7E6E6
$cs.Colors(5).RGB  = 13998939   # accent1  -> 5B9BD5
$cs.Colors(6).RGB  = 3243501    # accent2  -> ED7D31
$cs.Colors(7).RGB  = 10855845   # accent3  -> A5A5A5
$cs.Colors(8).RGB  = 49407      # accent4  -> FFC000
$cs.Colors(9).RGB  = 12874308   # accent5  -> 4472C4
$cs.Colors(10).RGB = 4697456    # accent6  -> 70AD47
$cs.Colors(11).RGB = 12673797   # hlink    -> 0563C1
$cs.Colors(12).RGB = 7491477    # folHlink -> 954F72
